$p = $ppt.ActivePresentation
$m = $p.Slides.Item(1).Master
foreach ($cl in $m.CustomLayouts) {
  Write-Output ("foreach layout shapes=" + $cl.Shapes.Count + " index=" + $cl.Index)
}
